# Friedman test sig addin: rename existing sheets and add 5 new
# "lnrN" result sheets (lnr3..lnr7), each holding a small linear
# regression summary table (names/coef/se/T/pval/r2/adj_r2/CI bounds).

function New-ResultSheet($wb, $headerSrc, $indexSrc, $name, $rows) {
    # rows: array of row-value arrays, each starting at column A
    #       (A = index, B = names, C..J = stats)
    $after = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $after)
    $ws.Name = $name

    $lastRow = $rows.Count + 1

    # Copy the bold/centered/bordered header style (B1:J1) and the
    # index-column style (A2:A{lastRow}) from an existing sheet.
    $headerSrc.Range("A1:J1").Copy()
    $ws.Range("A1:J1").PasteSpecial(-4122)

    $indexSrc.Range("A2:A$lastRow").Copy()
    $ws.Range("A2:A$lastRow").PasteSpecial(-4122)

    $headers = @("names", "coef", "se", "T", "pval", "r2", "adj_r2", "CI[2.5%]", "CI[97.5%]")
    for ($c = 0; $c -lt $headers.Count; $c++) {
        $ws.Cells.Item(1, $c + 2).Value = $headers[$c]
    }

    for ($r = 0; $r -lt $rows.Count; $r++) {
        $values = $rows[$r]
        for ($c = 0; $c -lt $values.Count; $c++) {
            $ws.Cells.Item($r + 2, $c + 1).Value = $values[$c]
        }
    }

    return $ws
}

$wb = $excel.ActiveWorkbook

# Rename the two existing sheets.
$wb.Worksheets.Item(1).Name = "lnr1"
$wb.Worksheets.Item(2).Name = "lnr2"

$headerSrc = $wb.Worksheets.Item(1)
$indexSrc = $wb.Worksheets.Item(1)

# lnr3
$rows3 = @(
    @(0, "Intercept", 0.5657302002800296, 0.07205185036140062, 7.851709531988657, [double]"3.070692274633801e-14", 0.8524422358406282, 0.8517820221307428, 0.4241277627072887, 0.7073326378527706),
    @(1, "Q4", 0.4837234431119523, 0.04449563947303484, 10.87125500028149, [double]"1.386884161679495e-24", 0.8524422358406282, 0.8517820221307428, 0.3962768198935847, 0.5711700663303199),
    @(2, "Q5", 0.387228343580688, 0.04258988160826051, 9.09202676688313, [double]"3.193423049236888e-18", 0.8524422358406282, 0.8517820221307428, 0.3035270781528823, 0.4709296090084938)
)
New-ResultSheet $wb $headerSrc $indexSrc "lnr3" $rows3

# lnr4
$rows4 = @(
    @(0, "Intercept", 0.6938496746467439, 0.07319332938513733, 9.479684562451908, [double]"1.510027432623504e-19", 0.8406773338610833, 0.8399644807687392, 0.5500039052041269, 0.8376954440893608),
    @(1, "Q9", 0.4393565877634387, 0.03694055258530152, 11.89361168187443, [double]"1.56821795995718e-28", 0.8406773338610833, 0.8399644807687392, 0.3667578652445098, 0.5119553102823675),
    @(2, "Q10", 0.4050326231213724, 0.0341976336593033, 11.84387864834575, [double]"2.46426757450648e-28", 0.8406773338610833, 0.8399644807687392, 0.3378245186648191, 0.4722407275779257)
)
New-ResultSheet $wb $headerSrc $indexSrc "lnr4" $rows4

# lnr5
$rows5 = @(
    @(0, "Intercept", 0.6197435897435907, 0.0780774296328991, 7.937551129147987, [double]"1.67135386548159e-14", 0.8251539689412495, 0.824763687621922, 0.4663001000578746, 0.7731870794293068),
    @(1, "Q4", 0.8570512820512818, 0.01863923975459291, 45.9810214008377, [double]"9.386016515947865e-172", 0.8251539689412495, 0.824763687621922, 0.8204200813767615, 0.8936824827258021)
)
New-ResultSheet $wb $headerSrc $indexSrc "lnr5" $rows5

# lnr6
$rows6 = @(
    @(0, "Intercept", 0.8028855565561189, 0.07713003279083254, 10.40950622610843, [double]"7.150524910775671e-23", 0.8134288261335085, 0.8130123726204137, 0.6513039606116213, 0.9544671525006165),
    @(1, "Q5", 0.8144943424298208, 0.01842941533422159, 44.19534356672649, [double]"1.949196879776997e-165", 0.8134288261335085, 0.8130123726204137, 0.7782755040914383, 0.8507131807682032)
)
New-ResultSheet $wb $headerSrc $indexSrc "lnr6" $rows6

# lnr7
$rows7 = @(
    @(0, "Intercept", 1.02299887659334, 0.07765880980104653, 13.17299195305918, [double]"1.004084975638781e-33", 0.7902579479064071, 0.789789773682984, 0.870378089296177, 1.175619663890502),
    @(1, "Q10", 0.763182336862776, 0.01857580415614014, 41.08475361000777, [double]"4.835855968827952e-154", 0.7902579479064071, 0.789789773682984, 0.7266758044777045, 0.7996888692478475)
)
New-ResultSheet $wb $headerSrc $indexSrc "lnr7" $rows7

Write-Host "Sheets now:"
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
